# "series to list. opening function."
# The "Job Title" column (e.g. "Python Coordinator") is split into a plain
# job "Title" column and a new "Programming Languages" column that lists
# out what used to be baked into the title text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "Job Title" column values --------------------------
# D1 header: "Job Title" -> "Title"
$ws.Range("D1").Value = "Title"
# D2 "General Manager" has no language, stays the same.
# D3 "Python Coordinator" -> "Coordinator" (the "Python" part moves to E3)
$ws.Range("D3").Value = "Coordinator"

# --- Add new "Programming Languages" column (E) --------------------------
$ws.Range("E1").Value = "Programming Languages"
$ws.Range("E2").Value = "Python, JavaScript"
$ws.Range("E3").Value = "Python"

# --- Header row formatting: bold font + taller row -----------------------
$ws.Range("A1:E1").Font.Bold = $true
$ws.Rows.Item(1).RowHeight = 28

# --- Update the on-screen selection to the new header row ----------------
$ws.Range("A1:E1").Select()
